# Update bitcoin_buys.xlsx after running on 2026-01-25
# Appends the newest weekly DCA buy as row 67 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 67

# Column A holds the purchase date as literal text (e.g. "01/18/2026" in the
# row above), not an Excel date serial. Prefixing with an apostrophe forces
# the COM layer to store it as text instead of auto-parsing it into a date.
$ws.Cells.Item($newRow, 1).Value = "'01/25/2026"

$ws.Cells.Item($newRow, 2).Value = 0.0005548899999999954
$ws.Cells.Item($newRow, 3).Value = 89206.86983005716
$ws.Cells.Item($newRow, 4).Value = 50
